# Add new power plant / electricity source rows to the two RQSD lookup
# sheets (RQSD-BRQSD and RQSD-RQSD), and clear the stray fill/font format
# that had been applied to the "see notes" cell on the About sheet.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBRQSD = $wb.Worksheets.Item("RQSD-BRQSD")
$wsRQSD  = $wb.Worksheets.Item("RQSD-RQSD")

# New electricity sources being added (issues #280 and #99).
$newSources = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

# Whether each new source qualifies under the (more permissive) BRQSD
# definition vs. the (stricter) RQSD definition.
$qualifiesBRQSD = @(1, 1, 1, 1, 1, 1)
$qualifiesRQSD  = @(0, 0, 0, 0, 0, 1)

$startRow = 19
for ($i = 0; $i -lt $newSources.Length; $i++) {
    $row = $startRow + $i

    $wsBRQSD.Cells.Item($row, 1).Value = $newSources[$i]
    $wsBRQSD.Cells.Item($row, 2).Value = $qualifiesBRQSD[$i]

    $wsRQSD.Cells.Item($row, 1).Value = $newSources[$i]
    $wsRQSD.Cells.Item($row, 2).Value = $qualifiesRQSD[$i]
}

# Move/restore the on-sheet selection so it lands just past the newly
# added rows, matching where the author's cursor ended up.
$wsBRQSD.Range("A25").Select()
$wsRQSD.Range("A35").Select()

# The "see notes" cell (B4) on the About sheet had picked up an extraneous
# font/fill style; clear its formatting back to the sheet default.
$wsAbout.Range("B4").ClearFormats()
$wsAbout.Activate()
